$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de) in row 2 show the status
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: column C (Status) in row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: column C (Status) in row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Shrink the Status columns to fit the new, shorter text ---
# "In Translation" is shorter than "Ready for handoff", so the report
# generator narrows the corresponding columns accordingly.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
